$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 409, shifting rows 409:470 down to 410:471
$ws.Rows("409:409").Insert()

# Populate the newly inserted row 409 with values
$ws.Cells.Item(409, 1).Value = 5
$ws.Cells.Item(409, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(409, 3).Value = "Maule"
$ws.Cells.Item(409, 4).Value = 45131
$ws.Cells.Item(409, 5).Value = 7
$ws.Cells.Item(409, 6).Value = "Fruta"
$ws.Cells.Item(409, 7).Value = 100101
$ws.Cells.Item(409, 8).Value = "Berries"
$ws.Cells.Item(409, 9).Value = 100101007
$ws.Cells.Item(409, 10).Value = "Kiwi"
$ws.Cells.Item(409, 11).Value = "Hayward"
$ws.Cells.Item(409, 12).Value = "Primera"
$ws.Cells.Item(409, 13).Value = 250
$ws.Cells.Item(409, 14).Value = 12000
$ws.Cells.Item(409, 15).Value = 12000
$ws.Cells.Item(409, 16).Value = 12000
$ws.Cells.Item(409, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(409, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(409, 19).Value = 667
$ws.Cells.Item(409, 20).Value = 18
